$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.189.17"
$ws.Range("E2").Value = "  -1.25%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.540.42"
$ws.Range("E3").Value = "  +0.38%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.16"
$ws.Range("E5").Value = "  +0.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.40"
$ws.Range("E6").Value = "  -2.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.539.80"
$ws.Range("E7").Value = "  +0.29%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.479"
$ws.Range("E9").Value = "  -0.25%  "

$ws.Range("E10").Value = "  -4.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "8.04"
$ws.Range("E11").Value = "  +2.76%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.414"
$ws.Range("E12").Value = "  -2.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.142.15"
$ws.Range("E13").Value = "  +0.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000208"
$ws.Range("E14").Value = "  -3.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.45"
$ws.Range("E15").Value = "  -4.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.539.20"
$ws.Range("E16").Value = "  +0.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.347.51"
$ws.Range("E17").Value = "  -1.16%  "

$ws.Range("E18").Value = "  -0.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.89"
$ws.Range("E19").Value = "  +0.39%  "

$ws.Range("E20").Value = "  -2.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.01"
$ws.Range("E21").Value = "  -2.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "425.78"
$ws.Range("E22").Value = "  -2.52%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.602"
$ws.Range("E23").Value = "  -1.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.59"
$ws.Range("E24").Value = "  -1.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.681.01"
$ws.Range("E25").Value = "  +0.59%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000120"
$ws.Range("E27").Value = "  -1.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.36"
$ws.Range("E28").Value = "  -5.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.01"
$ws.Range("E29").Value = "  -5.18%  "

$ws.Range("E30").Value = "  -1.26%  "

$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.48"
$ws.Range("E32").Value = "  -7.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.158"
$ws.Range("E33").Value = "  -6.76%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.36"
$ws.Range("E34").Value = "  -0.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.532.16"
$ws.Range("E35").Value = "  +0.66%  "

$ws.Range("E37").Value = "  -3.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.85"
$ws.Range("E38").Value = "  -2.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.61"
$ws.Range("E39").Value = "  -5.00%  "

$ws.Range("E40").Value = "  +0.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "169.49"
$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0862"
$ws.Range("E42").Value = "  -3.64%  "

$ws.Range("E43").Value = "  -4.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.893"
$ws.Range("E44").Value = "  -0.52%  "

$ws.Range("E45").Value = "  -9.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.28"
$ws.Range("E46").Value = "  -1.34%  "

$ws.Range("E47").Value = "  -8.16%  "

$ws.Range("E48").Value = "  -7.83%  "

$ws.Range("E49").Value = "  -1.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.19"
$ws.Range("E50").Value = "  -4.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.952"
$ws.Range("E51").Value = "  -4.41%  "
